# Update "想去人数" (want-to-go count, column F) figures to the latest
# scraped values across the four sheets of the workbook, matching the
# gh-pages data refresh generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 6046
$ws1.Range("F5").Value  = 77
$ws1.Range("F6").Value  = 66
$ws1.Range("F13").Value = 1625
$ws1.Range("F14").Value = 1625
$ws1.Range("F16").Value = 1669
$ws1.Range("F17").Value = 575
$ws1.Range("F20").Value = 4719
$ws1.Range("F21").Value = 119
$ws1.Range("F22").Value = 54
$ws1.Range("F25").Value = 828
$ws1.Range("F28").Value = 23
$ws1.Range("F29").Value = 2362
$ws1.Range("F40").Value = 1298
$ws1.Range("F41").Value = 1284

# --- Sheet 2: 演出 (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value  = 111
$ws2.Range("F13").Value = 112
$ws2.Range("F15").Value = 80

# --- Sheet 3: 本地生活 (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 777
$ws3.Range("F5").Value = 317

# --- Sheet 4: 全部类型 (All Types - combined listing) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 777
$ws4.Range("F7").Value  = 6046
$ws4.Range("F8").Value  = 77
$ws4.Range("F9").Value  = 66
$ws4.Range("F11").Value = 111
$ws4.Range("F22").Value = 1625
$ws4.Range("F24").Value = 1669
$ws4.Range("F25").Value = 80
$ws4.Range("F26").Value = 575
$ws4.Range("F29").Value = 4720
$ws4.Range("F30").Value = 54
$ws4.Range("F36").Value = 23
$ws4.Range("F37").Value = 2362
$ws4.Range("F49").Value = 1298
